$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C145").Value = 1.168885292845949
$ws.Range("C147").Value = 1.21486052238984
$ws.Range("C148").Value = 1.279036382914401
$ws.Range("C152").Value = 1.236072833604953
$ws.Range("C154").Value = 1.22415486547171
$ws.Range("C155").Value = 1.249077098369515
$ws.Range("C159").Value = 1.224516623136881
$ws.Range("C161").Value = 1.235914030567496
$ws.Range("C162").Value = 1.208582765258972
$ws.Range("C164").Value = 1.223339205264765
$ws.Range("C165").Value = 1.224129108051137
$ws.Range("C166").Value = 1.235573301768256
$ws.Range("C175").Value = 1.261125837412918
$ws.Range("C176").Value = 1.249623747082607
$ws.Range("C178").Value = 1.24948293342525
$ws.Range("C179").Value = 1.24989081645374
$ws.Range("C184").Value = 1.241279612420844
$ws.Range("C185").Value = 1.240343896002161
$ws.Range("C186").Value = 1.239565513799398
$ws.Range("C187").Value = 1.226372595977206
$ws.Range("C188").Value = 1.22690209228363
$ws.Range("C189").Value = 1.22779570041693
$ws.Range("C190").Value = 1.205641576628107
$ws.Range("C191").Value = 1.216840934226896
$ws.Range("C192").Value = 1.227545350993828
$ws.Range("C193").Value = 1.226995105162298
$ws.Range("C194").Value = 1.228583304476012
$ws.Range("C195").Value = 1.202035027553715
$ws.Range("C196").Value = 1.201453220846558
$ws.Range("C197").Value = 1.200773896767373
$ws.Range("C201").Value = 1.149225999358948
$ws.Range("C202").Value = 1.173906296295284
$ws.Range("C207").Value = 1.170979371306329
$ws.Range("C208").Value = 1.170410974266205
$ws.Range("C209").Value = 1.211399886366022
$ws.Range("C210").Value = 1.209062815901802
$ws.Range("C211").Value = 1.209767756797069
$ws.Range("C212").Value = 1.208126048465026
$ws.Range("C213").Value = 1.206143856762223
$ws.Range("C214").Value = 1.219685195887786
$ws.Range("C215").Value = 1.2294443467764
$ws.Range("C216").Value = 1.252166326270286
$ws.Range("C217").Value = 1.263445535777546
$ws.Range("C218").Value = 1.264352006479259
$ws.Range("C219").Value = 1.264027669521014
$ws.Range("C220").Value = 1.263709736548671
$ws.Range("C221").Value = 1.263976970690259
$ws.Range("C222").Value = 1.26387563390641
$ws.Range("C223").Value = 1.278942049283157
$ws.Range("C224").Value = 1.275572667932509
$ws.Range("C225").Value = 1.225329112322245
$ws.Range("C226").Value = 1.22508264222341
$ws.Range("C227").Value = 1.224625767606213
$ws.Range("C228").Value = 1.235961747676092
$ws.Range("C252").Value = 1.505024410882924
$ws.Range("C253").Value = 1.504056378694462
$ws.Range("C256").Value = 1.521293420386913
$ws.Range("C258").Value = 1.564888318697708
$ws.Range("C259").Value = 1.565272833530306
$ws.Range("C260").Value = 1.511848028050094
$ws.Range("C261").Value = 1.518667350800931
$ws.Range("C262").Value = 1.517943032579118
$ws.Range("C263").Value = 1.517641299295436
